$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab_7a_Quellen")

# Widen column B (target stored width 81.73828125 chars; ColumnWidth is
# quantized to an MDW-7 pixel grid by this engine, so 81 is the closest
# achievable input -> stored width 81.71428571428571).
$ws.Columns.Item(2).ColumnWidth = 81

# Row 4 (Q_AGEB): "AG Energiebilanzen" -> "AG Energiebilanzen e.V." (D/E)
$ws.Range("D4").Value = "AG Energiebilanzen e.V."
$ws.Range("E4").Value = "AG Energiebilanzen e.V."

# Row 37 (Q_JKI): reworded German institute names (B/D)
$ws.Range("B37").Value = "Institut für Pflanzenbau und Bodenkunde des Julius Kühn-Institut (JKI)"
$ws.Range("D37").Value = "Institut für Pflanzenbau und Bodenkunde des Julius Kühn-Institut"

# Row 48 (Q_UBALAWA): drop the "(LAWA)" suffix (D/E)
$ws.Range("D48").Value = "Umweltbundesamt nach Angaben der Bund/Länder Arbeitsgemeinschaft Wasser"
$ws.Range("E48").Value = "German Environment Agency on the basis of data from the German Working Group on Water Issues of the Länder and the Federal Government"

# Row 49 (Q_UG): comma -> "der" (B/D)
$ws.Range("B49").Value = "Institut für Landschaftsökologie und Ressourcenmanagement der Justus-Liebig-Universität Gießen"
$ws.Range("D49").Value = "Institut für Landschaftsökologie und Ressourcenmanagement der Justus-Liebig-Universität Gießen"
